# =====================================================================
# PlayerPerformance_3935.xlsx -- additional scraping update
#   1. Add a new "Player Info" sheet (first tab) with player bio data.
#   2. Rewrite "ODI Batting"  MATCH_CARD_LINK -> MATCH_CODE (bare code).
#   3. Rewrite "ODI Bowling"  MATCH_CARD_LINK -> MATCH_CODE (bare code).
#   4. Add a new "ODI Batting Extra" sheet (last tab) with additional
#      per-innings batting detail.
# =====================================================================

$wb = $excel.ActiveWorkbook

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------
# 1. "Player Info" sheet -- inserted before "ODI Batting" (becomes tab 1)
# ---------------------------------------------------------------------
$wsPlayerInfo = $wb.Worksheets.Add($wsBatting)
$wsPlayerInfo.Name = "Player Info"

$wsPlayerInfo.Range("A1").Value = "ID"
$wsPlayerInfo.Range("B1").Value = "NAME"
$wsPlayerInfo.Range("C1").Value = "BATTING_HAND"
$wsPlayerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRng = $wsPlayerInfo.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

$wsPlayerInfo.Range("A2").NumberFormat = "@"
$wsPlayerInfo.Range("A2").Value = "3935"
$wsPlayerInfo.Range("B2").Value = "Amir Hamza"
$wsPlayerInfo.Range("C2").Value = "Right Handed"
$wsPlayerInfo.Range("D2").Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------
# 2. "ODI Batting" -- MATCH_CARD_LINK column becomes MATCH_CODE
# ---------------------------------------------------------------------
$wsBatting.Range("D1").Value = "MATCH_CODE"

for ($r = 2; $r -le 32; $r++) {
    $linkCell = $wsBatting.Cells.Item($r, 4)
    $url = $linkCell.Text
    if ($url -ne "") {
        $code = $url -replace '.*MatchCode=', ''
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $code
    }

    # Drop the INNING_NUMBER cells that only ever held an empty string --
    # they should no longer be present at all.
    $inningCell = $wsBatting.Cells.Item($r, 2)
    if ($inningCell.Text -eq "") {
        $inningCell.Value = ""
    }
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling" -- MATCH_CARD_LINK column becomes MATCH_CODE
# ---------------------------------------------------------------------
$wsBowling.Range("B1").Value = "MATCH_CODE"

for ($r = 2; $r -le 32; $r++) {
    $linkCell = $wsBowling.Cells.Item($r, 2)
    $url = $linkCell.Text
    if ($url -ne "") {
        $code = $url -replace '.*MatchCode=', ''
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $code
    }
}

# ---------------------------------------------------------------------
# 4. "ODI Batting Extra" -- inserted after "ODI Bowling" (last tab)
# ---------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Add($null, $wsBowling)
$wsExtra.Name = "ODI Batting Extra"

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

$extraHeaderRng = $wsExtra.Range("A1:F1")
$extraHeaderRng.Font.Bold = $true
$extraHeaderRng.HorizontalAlignment = -4108
$extraHeaderRng.VerticalAlignment = -4160
$extraHeaderRng.Borders.LineStyle = 1

$extraData = @(
    @("3843", "8", "0", "0", "3.28%", "NO"),
    @("3846", "11", "", "", "", "NO"),
    @("3847", "", "", "", "", "NO"),
    @("3849", "", "", "", "", "NO"),
    @("3850", "11", "0", "0", "", "NO"),
    @("3864", "10", "0", "0", "1.53%", "YES"),
    @("3867", "10", "", "", "", "NO"),
    @("3869", "", "", "", "", "NO"),
    @("3871", "", "", "", "", "NO"),
    @("3873", "11", "0", "0", "", "NO"),
    @("3991", "10", "0", "0", "1.59%", "NO"),
    @("3993", "11", "0", "0", "", "NO"),
    @("3994", "10", "", "", "", "NO"),
    @("3998", "11", "0", "0", "", "NO"),
    @("4006", "", "", "", "", "NO"),
    @("4007", "11", "0", "1", "3.18%", "NO"),
    @("4008", "10", "", "", "", "NO"),
    @("4009", "11", "0", "0", "", "NO"),
    @("4040", "", "", "", "", ""),
    @("4043", "", "", "", "", "")
)

$r = 2
foreach ($row in $extraData) {
    $wsExtra.Cells.Item($r, 1).NumberFormat = "@"
    $wsExtra.Cells.Item($r, 1).Value = $row[0]

    if ($row[1] -ne "") {
        $wsExtra.Cells.Item($r, 2).Value = [int]$row[1]
    }

    $wsExtra.Cells.Item($r, 3).NumberFormat = "@"
    $wsExtra.Cells.Item($r, 3).Value = $row[2]

    $wsExtra.Cells.Item($r, 4).NumberFormat = "@"
    $wsExtra.Cells.Item($r, 4).Value = $row[3]

    $wsExtra.Cells.Item($r, 5).NumberFormat = "@"
    $wsExtra.Cells.Item($r, 5).Value = $row[4]

    $wsExtra.Cells.Item($r, 6).NumberFormat = "@"
    $wsExtra.Cells.Item($r, 6).Value = $row[5]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Restore the originally-active sheet/selection.
# ---------------------------------------------------------------------
$wsPlayerInfo.Activate()
$wsPlayerInfo.Range("A1").Select()
